# Auto-generated Word COM-interop script to update the daily math worksheet
# for two-digit x two-digit multiplication: refreshes the date heading and
# all 25 multiplication problems/answers in the table.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-09-30 Monday" "2024-10-01 Tuesday"
Replace-Text "57×93=5301" "66×29=1914"
Replace-Text "82×25=2050" "58×96=5568"
Replace-Text "53×92=4876" "46×75=3450"
Replace-Text "75×86=6450" "43×84=3612"
Replace-Text "11×22=242" "17×37=629"
Replace-Text "46×85=3910" "48×53=2544"
Replace-Text "89×45=4005" "65×74=4810"
Replace-Text "43×18=774" "27×19=513"
Replace-Text "63×63=3969" "43×38=1634"
Replace-Text "43×24=1032" "36×42=1512"
Replace-Text "64×53=3392" "35×18=630"
Replace-Text "18×70=1260" "86×31=2666"
Replace-Text "63×59=3717" "54×63=3402"
Replace-Text "63×62=3906" "12×93=1116"
Replace-Text "69×76=5244" "66×50=3300"
Replace-Text "53×45=2385" "63×67=4221"
Replace-Text "86×53=4558" "62×66=4092"
Replace-Text "24×60=1440" "88×41=3608"
Replace-Text "97×61=5917" "83×49=4067"
Replace-Text "93×30=2790" "13×34=442"
Replace-Text "65×98=6370" "60×63=3780"
Replace-Text "45×51=2295" "17×77=1309"
Replace-Text "38×84=3192" "13×48=624"
Replace-Text "37×72=2664" "24×48=1152"
Replace-Text "12×84=1008" "17×51=867"

Write-Output "Done: updated date heading and 25 multiplication cells."
